$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row33 = New-Object "object[,]" 1,50
$row33[0,0] = 8
$row33[0,1] = 9
$row33[0,2] = 10
$row33[0,3] = 11
$row33[0,4] = 12
$row33[0,5] = 13
$row33[0,6] = 14
$row33[0,7] = 15
$row33[0,8] = 16
$row33[0,9] = 17
$row33[0,10] = 18
$row33[0,11] = 19
$row33[0,12] = 20
$row33[0,13] = 21
$row33[0,14] = 22
$row33[0,15] = 23
$row33[0,16] = 24
$row33[0,17] = 25
$row33[0,18] = 26
$row33[0,19] = 27
$row33[0,20] = 28
$row33[0,21] = 29
$row33[0,22] = 30
$row33[0,23] = 31
$row33[0,24] = 32
$row33[0,25] = 33
$row33[0,26] = 34
$row33[0,27] = 35
$row33[0,28] = 36
$row33[0,29] = 37
$row33[0,30] = 38
$row33[0,31] = 39
$row33[0,32] = 40
$row33[0,33] = 41
$row33[0,34] = 42
$row33[0,35] = 43
$row33[0,36] = 44
$row33[0,37] = 45
$row33[0,38] = 46
$row33[0,39] = 47
$row33[0,40] = 48
$row33[0,41] = 49
$row33[0,42] = 50
$row33[0,43] = 51
$row33[0,44] = 52
$row33[0,45] = 53
$row33[0,46] = 54
$row33[0,47] = 55
$row33[0,48] = 56
$row33[0,49] = 57
$ws.Range("A33:AX33").Value = $row33

$row34 = New-Object "object[,]" 1,50
$row34[0,0] = 1.1428571428571428
$row34[0,1] = 1.2857142857142856
$row34[0,2] = 1.4285714285714286
$row34[0,3] = 1.5714285714285714
$row34[0,4] = 1.7142857142857144
$row34[0,5] = 1.8571428571428572
$row34[0,6] = 2
$row34[0,7] = 2.1428571428571428
$row34[0,8] = 2.2857142857142856
$row34[0,9] = 2.4285714285714288
$row34[0,10] = 2.5714285714285712
$row34[0,11] = 2.7142857142857144
$row34[0,12] = 2.8571428571428572
$row34[0,13] = 3
$row34[0,14] = 3.1428571428571428
$row34[0,15] = 3.2857142857142856
$row34[0,16] = 3.4285714285714284
$row34[0,17] = 3.5714285714285716
$row34[0,18] = 3.7142857142857144
$row34[0,19] = 3.8571428571428572
$row34[0,20] = 4
$row34[0,21] = 4.1428571428571423
$row34[0,22] = 4.2857142857142856
$row34[0,23] = 4.4285714285714288
$row34[0,24] = 4.5714285714285712
$row34[0,25] = 4.7142857142857144
$row34[0,26] = 4.8571428571428577
$row34[0,27] = 5
$row34[0,28] = 5.1428571428571432
$row34[0,29] = 5.2857142857142856
$row34[0,30] = 5.4285714285714288
$row34[0,31] = 5.5714285714285712
$row34[0,32] = 5.7142857142857144
$row34[0,33] = 5.8571428571428568
$row34[0,34] = 6
$row34[0,35] = 6.1428571428571432
$row34[0,36] = 6.2857142857142856
$row34[0,37] = 6.4285714285714288
$row34[0,38] = 6.5714285714285712
$row34[0,39] = 6.7142857142857144
$row34[0,40] = 6.8571428571428568
$row34[0,41] = 7
$row34[0,42] = 7.1428571428571432
$row34[0,43] = 7.2857142857142856
$row34[0,44] = 7.4285714285714288
$row34[0,45] = 7.5714285714285712
$row34[0,46] = 7.7142857142857144
$row34[0,47] = 7.8571428571428568
$row34[0,48] = 8
$row34[0,49] = 8.1428571428571423
$ws.Range("A34:AX34").Value = $row34

$row35 = New-Object "object[,]" 1,50
$row35[0,0] = 6.8550000000000004
$row35[0,1] = 17.55
$row35[0,2] = 24.427
$row35[0,3] = 31.541
$row35[0,4] = 34.264000000000003
$row35[0,5] = 36.984999999999999
$row35[0,6] = 25.550999999999998
$row35[0,7] = 38.777999999999999
$row35[0,8] = 38.061
$row35[0,9] = 37.593000000000004
$row35[0,10] = 37.298999999999999
$row35[0,11] = 37.563000000000002
$row35[0,12] = 36.406999999999996
$row35[0,13] = 37.578000000000003
$row35[0,14] = 37.390999999999998
$row35[0,15] = 36.863
$row35[0,16] = 36.045000000000002
$row35[0,17] = 34.987000000000002
$row35[0,18] = 33.74
$row35[0,19] = 32.351999999999997
$row35[0,20] = 29.940999999999999
$row35[0,21] = 31.099
$row35[0,22] = 28.72
$row35[0,23] = 24.3
$row35[0,24] = 21.920999999999999
$row35[0,25] = 25.283999999999999
$row35[0,26] = 28.646999999999998
$row35[0,27] = 28.439
$row35[0,28] = 27.9
$row35[0,29] = 27.158999999999999
$row35[0,30] = 26.344999999999999
$row35[0,31] = 25.042999999999999
$row35[0,32] = 23.23
$row35[0,33] = 21.619
$row35[0,34] = 20.922999999999998
$row35[0,35] = 23.452000000000002
$row35[0,36] = 23.355
$row35[0,37] = 23.088999999999999
$row35[0,38] = 22.693999999999999
$row35[0,39] = 22.209
$row35[0,40] = 21.672999999999998
$row35[0,41] = 21.125
$row35[0,42] = 20.605
$row35[0,43] = 20.027999999999999
$row35[0,44] = 19.36
$row35[0,45] = 18.713999999999999
$row35[0,46] = 18.202999999999999
$row35[0,47] = 17.821000000000002
$row35[0,48] = 17.5
$row35[0,49] = 17.260999999999999
$ws.Range("A35:AX35").Value = $row35

$row38 = New-Object "object[,]" 1,50
$row38[0,0] = 9
$row38[0,1] = 10
$row38[0,2] = 11
$row38[0,3] = 12
$row38[0,4] = 13
$row38[0,5] = 14
$row38[0,6] = 15
$row38[0,7] = 16
$row38[0,8] = 17
$row38[0,9] = 18
$row38[0,10] = 19
$row38[0,11] = 20
$row38[0,12] = 21
$row38[0,13] = 22
$row38[0,14] = 23
$row38[0,15] = 24
$row38[0,16] = 25
$row38[0,17] = 26
$row38[0,18] = 27
$row38[0,19] = 28
$row38[0,20] = 29
$row38[0,21] = 30
$row38[0,22] = 31
$row38[0,23] = 32
$row38[0,24] = 33
$row38[0,25] = 34
$row38[0,26] = 35
$row38[0,27] = 36
$row38[0,28] = 37
$row38[0,29] = 38
$row38[0,30] = 39
$row38[0,31] = 40
$row38[0,32] = 41
$row38[0,33] = 42
$row38[0,34] = 43
$row38[0,35] = 44
$row38[0,36] = 45
$row38[0,37] = 46
$row38[0,38] = 47
$row38[0,39] = 48
$row38[0,40] = 49
$row38[0,41] = 50
$row38[0,42] = 51
$row38[0,43] = 52
$row38[0,44] = 53
$row38[0,45] = 54
$row38[0,46] = 55
$row38[0,47] = 56
$row38[0,48] = 57
$row38[0,49] = 58
$ws.Range("A38:AX38").Value = $row38

$row39 = New-Object "object[,]" 1,50
$row39[0,0] = 1.125
$row39[0,1] = 1.25
$row39[0,2] = 1.375
$row39[0,3] = 1.5
$row39[0,4] = 1.625
$row39[0,5] = 1.75
$row39[0,6] = 1.875
$row39[0,7] = 2
$row39[0,8] = 2.125
$row39[0,9] = 2.25
$row39[0,10] = 2.375
$row39[0,11] = 2.5
$row39[0,12] = 2.625
$row39[0,13] = 2.75
$row39[0,14] = 2.875
$row39[0,15] = 3
$row39[0,16] = 3.125
$row39[0,17] = 3.25
$row39[0,18] = 3.375
$row39[0,19] = 3.5
$row39[0,20] = 3.625
$row39[0,21] = 3.75
$row39[0,22] = 3.875
$row39[0,23] = 4
$row39[0,24] = 4.125
$row39[0,25] = 4.25
$row39[0,26] = 4.375
$row39[0,27] = 4.5
$row39[0,28] = 4.625
$row39[0,29] = 4.75
$row39[0,30] = 4.875
$row39[0,31] = 5
$row39[0,32] = 5.125
$row39[0,33] = 5.25
$row39[0,34] = 5.375
$row39[0,35] = 5.5
$row39[0,36] = 5.625
$row39[0,37] = 5.75
$row39[0,38] = 5.875
$row39[0,39] = 6
$row39[0,40] = 6.125
$row39[0,41] = 6.25
$row39[0,42] = 6.375
$row39[0,43] = 6.5
$row39[0,44] = 6.625
$row39[0,45] = 6.75
$row39[0,46] = 6.875
$row39[0,47] = 7
$row39[0,48] = 7.125
$row39[0,49] = 7.25
$ws.Range("A39:AX39").Value = $row39

$row40 = New-Object "object[,]" 1,50
$row40[0,0] = 6.8819999999999997
$row40[0,1] = 18.625
$row40[0,2] = 26.263999999999999
$row40[0,3] = 34.796999999999997
$row40[0,4] = 35.384
$row40[0,5] = 36.887999999999998
$row40[0,6] = 37.512
$row40[0,7] = 28.463000000000001
$row40[0,8] = 38.103999999999999
$row40[0,9] = 37.331000000000003
$row40[0,10] = 36.798999999999999
$row40[0,11] = 37.634
$row40[0,12] = 36.225999999999999
$row40[0,13] = 35.076999999999998
$row40[0,14] = 34.304000000000002
$row40[0,15] = 33.649000000000001
$row40[0,16] = 33.079000000000001
$row40[0,17] = 32.561
$row40[0,18] = 32.061999999999998
$row40[0,19] = 31.547000000000001
$row40[0,20] = 30.984000000000002
$row40[0,21] = 30.373000000000001
$row40[0,22] = 29.745999999999999
$row40[0,23] = 29.113
$row40[0,24] = 28.489000000000001
$row40[0,25] = 27.884
$row40[0,26] = 27.31
$row40[0,27] = 26.797999999999998
$row40[0,28] = 26.350999999999999
$row40[0,29] = 25.937000000000001
$row40[0,30] = 25.526
$row40[0,31] = 25.087
$row40[0,32] = 24.59
$row40[0,33] = 24.004999999999999
$row40[0,34] = 23.3
$row40[0,35] = 22.16
$row40[0,36] = 21.483000000000001
$row40[0,37] = 20.917999999999999
$row40[0,38] = 20.440000000000001
$row40[0,39] = 20.023
$row40[0,40] = 19.64
$row40[0,41] = 19.265000000000001
$row40[0,42] = 18.870999999999999
$row40[0,43] = 18.465
$row40[0,44] = 18.071000000000002
$row40[0,45] = 17.690000000000001
$row40[0,46] = 17.321999999999999
$row40[0,47] = 16.969000000000001
$row40[0,48] = 16.63
$row40[0,49] = 16.308
$ws.Range("A40:AX40").Value = $row40

$row43 = New-Object "object[,]" 1,50
$row43[0,0] = 10
$row43[0,1] = 11
$row43[0,2] = 12
$row43[0,3] = 13
$row43[0,4] = 14
$row43[0,5] = 15
$row43[0,6] = 16
$row43[0,7] = 17
$row43[0,8] = 18
$row43[0,9] = 19
$row43[0,10] = 20
$row43[0,11] = 21
$row43[0,12] = 22
$row43[0,13] = 23
$row43[0,14] = 24
$row43[0,15] = 25
$row43[0,16] = 26
$row43[0,17] = 27
$row43[0,18] = 28
$row43[0,19] = 29
$row43[0,20] = 30
$row43[0,21] = 31
$row43[0,22] = 32
$row43[0,23] = 33
$row43[0,24] = 34
$row43[0,25] = 35
$row43[0,26] = 36
$row43[0,27] = 37
$row43[0,28] = 38
$row43[0,29] = 39
$row43[0,30] = 40
$row43[0,31] = 41
$row43[0,32] = 42
$row43[0,33] = 43
$row43[0,34] = 44
$row43[0,35] = 45
$row43[0,36] = 46
$row43[0,37] = 47
$row43[0,38] = 48
$row43[0,39] = 49
$row43[0,40] = 50
$row43[0,41] = 51
$row43[0,42] = 52
$row43[0,43] = 53
$row43[0,44] = 54
$row43[0,45] = 55
$row43[0,46] = 56
$row43[0,47] = 57
$row43[0,48] = 58
$row43[0,49] = 59
$ws.Range("A43:AX43").Value = $row43

$row44 = New-Object "object[,]" 1,50
$row44[0,0] = 1.1111111111111112
$row44[0,1] = 1.2222222222222223
$row44[0,2] = 1.3333333333333333
$row44[0,3] = 1.4444444444444444
$row44[0,4] = 1.5555555555555556
$row44[0,5] = 1.6666666666666665
$row44[0,6] = 1.7777777777777777
$row44[0,7] = 1.8888888888888888
$row44[0,8] = 2
$row44[0,9] = 2.1111111111111112
$row44[0,10] = 2.2222222222222223
$row44[0,11] = 2.333333333333333
$row44[0,12] = 2.4444444444444446
$row44[0,13] = 2.5555555555555554
$row44[0,14] = 2.666666666666667
$row44[0,15] = 2.7777777777777777
$row44[0,16] = 2.8888888888888888
$row44[0,17] = 3
$row44[0,18] = 3.1111111111111112
$row44[0,19] = 3.2222222222222223
$row44[0,20] = 3.3333333333333335
$row44[0,21] = 3.4444444444444446
$row44[0,22] = 3.5555555555555554
$row44[0,23] = 3.6666666666666665
$row44[0,24] = 3.7777777777777777
$row44[0,25] = 3.8888888888888888
$row44[0,26] = 4
$row44[0,27] = 4.1111111111111107
$row44[0,28] = 4.2222222222222223
$row44[0,29] = 4.3333333333333339
$row44[0,30] = 4.4444444444444446
$row44[0,31] = 4.5555555555555554
$row44[0,32] = 4.6666666666666661
$row44[0,33] = 4.7777777777777777
$row44[0,34] = 4.8888888888888893
$row44[0,35] = 5
$row44[0,36] = 5.1111111111111107
$row44[0,37] = 5.2222222222222223
$row44[0,38] = 5.333333333333333
$row44[0,39] = 5.4444444444444446
$row44[0,40] = 5.5555555555555554
$row44[0,41] = 5.666666666666667
$row44[0,42] = 5.7777777777777777
$row44[0,43] = 5.8888888888888893
$row44[0,44] = 6
$row44[0,45] = 6.1111111111111107
$row44[0,46] = 6.2222222222222223
$row44[0,47] = 6.333333333333333
$row44[0,48] = 6.4444444444444446
$row44[0,49] = 6.5555555555555554
$ws.Range("A44:AX44").Value = $row44

$row45 = New-Object "object[,]" 1,50
$row45[0,0] = 7.024
$row45[0,1] = 17.890999999999998
$row45[0,2] = 20.375
$row45[0,3] = 31.635000000000002
$row45[0,4] = 34.787999999999997
$row45[0,5] = 36.567
$row45[0,6] = 37.472000000000001
$row45[0,7] = 37.548000000000002
$row45[0,8] = 30.902999999999999
$row45[0,9] = 37.173000000000002
$row45[0,10] = 36.15
$row45[0,11] = 36.21
$row45[0,12] = 34.792000000000002
$row45[0,13] = 34.313000000000002
$row45[0,14] = 33.750999999999998
$row45[0,15] = 33.134999999999998
$row45[0,16] = 32.481999999999999
$row45[0,17] = 31.809000000000001
$row45[0,18] = 31.132999999999999
$row45[0,19] = 30.471
$row45[0,20] = 29.841000000000001
$row45[0,21] = 29.234000000000002
$row45[0,22] = 28.632999999999999
$row45[0,23] = 28.041
$row45[0,24] = 27.46
$row45[0,25] = 26.895
$row45[0,26] = 26.349
$row45[0,27] = 25.841999999999999
$row45[0,28] = 25.379000000000001
$row45[0,29] = 24.940999999999999
$row45[0,30] = 24.506
$row45[0,31] = 24.056000000000001
$row45[0,32] = 23.568999999999999
$row45[0,33] = 23.027000000000001
$row45[0,34] = 22.408000000000001
$row45[0,35] = 21.574999999999999
$row45[0,36] = 21.05
$row45[0,37] = 20.632999999999999
$row45[0,38] = 20.268000000000001
$row45[0,39] = 19.899000000000001
$row45[0,40] = 19.535
$row45[0,41] = 19.187000000000001
$row45[0,42] = 18.800999999999998
$row45[0,43] = 18.311
$row45[0,44] = 17.77
$row45[0,45] = 17.302
$row45[0,46] = 16.925000000000001
$row45[0,47] = 16.581
$row45[0,48] = 16.283000000000001
$row45[0,49] = 16.042000000000002
$ws.Range("A45:AX45").Value = $row45

$row48 = New-Object "object[,]" 1,50
$row48[0,0] = 11
$row48[0,1] = 12
$row48[0,2] = 13
$row48[0,3] = 14
$row48[0,4] = 15
$row48[0,5] = 16
$row48[0,6] = 17
$row48[0,7] = 18
$row48[0,8] = 19
$row48[0,9] = 20
$row48[0,10] = 21
$row48[0,11] = 22
$row48[0,12] = 23
$row48[0,13] = 24
$row48[0,14] = 25
$row48[0,15] = 26
$row48[0,16] = 27
$row48[0,17] = 28
$row48[0,18] = 29
$row48[0,19] = 30
$row48[0,20] = 31
$row48[0,21] = 32
$row48[0,22] = 33
$row48[0,23] = 34
$row48[0,24] = 35
$row48[0,25] = 36
$row48[0,26] = 37
$row48[0,27] = 38
$row48[0,28] = 39
$row48[0,29] = 40
$row48[0,30] = 41
$row48[0,31] = 42
$row48[0,32] = 43
$row48[0,33] = 44
$row48[0,34] = 45
$row48[0,35] = 46
$row48[0,36] = 47
$row48[0,37] = 48
$row48[0,38] = 49
$row48[0,39] = 50
$row48[0,40] = 51
$row48[0,41] = 52
$row48[0,42] = 53
$row48[0,43] = 54
$row48[0,44] = 55
$row48[0,45] = 56
$row48[0,46] = 57
$row48[0,47] = 58
$row48[0,48] = 59
$row48[0,49] = 60
$ws.Range("A48:AX48").Value = $row48

$row49 = New-Object "object[,]" 1,50
$row49[0,0] = 1.1000000000000001
$row49[0,1] = 1.2
$row49[0,2] = 1.3
$row49[0,3] = 1.4
$row49[0,4] = 1.5
$row49[0,5] = 1.6
$row49[0,6] = 1.7
$row49[0,7] = 1.8
$row49[0,8] = 1.9
$row49[0,9] = 2
$row49[0,10] = 2.1
$row49[0,11] = 2.2000000000000002
$row49[0,12] = 2.2999999999999998
$row49[0,13] = 2.4
$row49[0,14] = 2.5
$row49[0,15] = 2.6
$row49[0,16] = 2.7
$row49[0,17] = 2.8
$row49[0,18] = 2.9
$row49[0,19] = 3
$row49[0,20] = 3.1
$row49[0,21] = 3.2
$row49[0,22] = 3.3
$row49[0,23] = 3.4
$row49[0,24] = 3.5
$row49[0,25] = 3.6
$row49[0,26] = 3.7
$row49[0,27] = 3.8
$row49[0,28] = 3.9
$row49[0,29] = 4
$row49[0,30] = 4.0999999999999996
$row49[0,31] = 4.2
$row49[0,32] = 4.3
$row49[0,33] = 4.4000000000000004
$row49[0,34] = 4.5
$row49[0,35] = 4.5999999999999996
$row49[0,36] = 4.7
$row49[0,37] = 4.8
$row49[0,38] = 4.9000000000000004
$row49[0,39] = 5
$row49[0,40] = 5.0999999999999996
$row49[0,41] = 5.2
$row49[0,42] = 5.3
$row49[0,43] = 5.4
$row49[0,44] = 5.5
$row49[0,45] = 5.6
$row49[0,46] = 5.7
$row49[0,47] = 5.8
$row49[0,48] = 5.9
$row49[0,49] = 6
$ws.Range("A49:AX49").Value = $row49

$row50 = New-Object "object[,]" 1,50
$row50[0,0] = 7.1959999999999997
$row50[0,1] = 15.981999999999999
$row50[0,2] = 26.283999999999999
$row50[0,3] = 32.046999999999997
$row50[0,4] = 37.659999999999997
$row50[0,5] = 34.978999999999999
$row50[0,6] = 36.915999999999997
$row50[0,7] = 37.012
$row50[0,8] = 36.354999999999997
$row50[0,9] = 32.981000000000002
$row50[0,10] = 35.247999999999998
$row50[0,11] = 35.186999999999998
$row50[0,12] = 34.131
$row50[0,13] = 33.435000000000002
$row50[0,14] = 32.741
$row50[0,15] = 32.042999999999999
$row50[0,16] = 31.344999999999999
$row50[0,17] = 30.654
$row50[0,18] = 29.975000000000001
$row50[0,19] = 29.312999999999999
$row50[0,20] = 28.673999999999999
$row50[0,21] = 28.053999999999998
$row50[0,22] = 27.446000000000002
$row50[0,23] = 26.853000000000002
$row50[0,24] = 26.277000000000001
$row50[0,25] = 25.72
$row50[0,26] = 25.184000000000001
$row50[0,27] = 24.696999999999999
$row50[0,28] = 24.265999999999998
$row50[0,29] = 23.866
$row50[0,30] = 23.47
$row50[0,31] = 23.050999999999998
$row50[0,32] = 22.584
$row50[0,33] = 22.042000000000002
$row50[0,34] = 21.398
$row50[0,35] = 20.413
$row50[0,36] = 19.763999999999999
$row50[0,37] = 19.218
$row50[0,38] = 18.783000000000001
$row50[0,39] = 18.466000000000001
$row50[0,40] = 18.279
$row50[0,41] = 18.143000000000001
$row50[0,42] = 17.931000000000001
$row50[0,43] = 17.300999999999998
$row50[0,44] = 16.332999999999998
$row50[0,45] = 15.535
$row50[0,46] = 15.010999999999999
$row50[0,47] = 14.547000000000001
$row50[0,48] = 14.176
$row50[0,49] = 13.928000000000001
$ws.Range("A50:AX50").Value = $row50

$ws.Activate()
$ws.Range("C54").Select()
